$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking price strings that must remain text (matching the
# source data which used western "thousands dot" / truncated-decimal formatting).
# Force text number-format before assigning so Excel does not coerce them to numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.849.58"
$ws.Range("E2").Value = "  -1.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.891.69"
$ws.Range("E3").Value = "  -1.15%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7579"
$ws.Range("E5").Value = "  +2.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "239.55"
$ws.Range("E6").Value = "  -1.76%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9998"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3037"
$ws.Range("E8").Value = "  -3.56%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.41"
$ws.Range("E9").Value = "  -6.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06806"
$ws.Range("E10").Value = "  -2.88%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07978"
$ws.Range("E11").Value = "  +0.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7472"
$ws.Range("E12").Value = "  -4.70%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.895.77"
$ws.Range("E13").Value = "  -1.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.200"
$ws.Range("E14").Value = "  -1.90%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "90.90"
$ws.Range("E15").Value = "  -1.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.856.81"
$ws.Range("E16").Value = "  -1.32%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.88"
$ws.Range("E17").Value = "  -3.36%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.953"
$ws.Range("E18").Value = "  +1.48%  "
$ws.Range("E19").Value = "  -2.43%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "234.71"
$ws.Range("E20").Value = "  -4.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9993"
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.147.10"
$ws.Range("E22").Value = "  -2.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9993"
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.945"
$ws.Range("E24").Value = "  +4.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.217"
$ws.Range("E25").Value = "  -2.87%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.46"
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("E27").Value = "  -1.96%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1300"
$ws.Range("E28").Value = "  +1.77%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.032"
$ws.Range("E29").Value = "  -4.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.342"
$ws.Range("E30").Value = "  -1.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.514"
$ws.Range("E31").Value = "  -2.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.272"
$ws.Range("E32").Value = "  -1.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.019"
$ws.Range("E33").Value = "  -1.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05334"
$ws.Range("E34").Value = "  +1.83%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.242"
$ws.Range("E35").Value = "  -5.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7264"
$ws.Range("E36").Value = "  -3.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.722"
$ws.Range("E37").Value = "  -1.41%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01923"
$ws.Range("E38").Value = "  -1.35%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.775"
$ws.Range("E39").Value = "  -0.78%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.162"
$ws.Range("E40").Value = "  -4.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4396"
$ws.Range("E41").Value = "  -2.58%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "71.96"
$ws.Range("E42").Value = "  -5.70%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.911"
$ws.Range("E43").Value = "  -2.47%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9994"
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8244"
$ws.Range("E45").Value = "  -1.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.01"
$ws.Range("E46").Value = "  -0.37%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.554"
$ws.Range("E47").Value = "  -3.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.749"
$ws.Range("E48").Value = "  -1.57%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.042.49"
$ws.Range("E49").Value = "  -1.63%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.08"
$ws.Range("E50").Value = "  -3.10%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05957"
$ws.Range("E51").Value = "  -0.84%  "
